$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4 - new data row (quote-prefixed text, mirroring the existing rows' formatting)
$ws.Range("A4").Value = "'3"
$ws.Range("B4").Value = "'95400152"
$ws.Range("B4").Font.Name = "Calibri"
$ws.Range("C4").Value = "'1"
$ws.Range("D4").Value = "'sandrita69"
$ws.Range("D4").Font.Name = "Calibri"
$ws.Range("E4").Value = "'1234"
$ws.Range("E4").Font.Name = "Calibri"
$ws.Range("F4").Value = "'4321"
$ws.Range("F4").Font.Name = "Calibri"
$ws.Range("G4").Value = "'Acierto"
$ws.Range("H4").Value = "'001"
$ws.Range("I4").Value = "'0370"
$ws.Range("J4").Value = "'NO ERROR"
$ws.Range("K4").Value = "'bolp"
$ws.Range("L4").Value = "'ACTIVO"

# Column L got resized (best-fit) once the new data pushed its effective width
$ws.Columns.Item(12).AutoFit()

# Cursor ends up parked on K1 after the edit
$ws.Range("K1").Select()
